$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'257.02"
$ws.Range("E2").Value = "'0.60%"
$ws.Range("D3").Value = "'27.04"
$ws.Range("E3").Value = "'-3.59%"
$ws.Range("D4").Value = "'4.629"
$ws.Range("E4").Value = "'-11.15%"
$ws.Range("D5").Value = "'0.05895"
$ws.Range("E5").Value = "'0.65%"
$ws.Range("D6").Value = "'6.642"
$ws.Range("E6").Value = "'-0.75%"
$ws.Range("E7").Value = "'-0.07%"
$ws.Range("D8").Value = "'0.9408"
$ws.Range("E8").Value = "'-1.94%"
$ws.Range("D9").Value = "'0.1404"
$ws.Range("E9").Value = "'-0.27%"
$ws.Range("D10").Value = "'0.03878"
$ws.Range("E10").Value = "'11.11%"
$ws.Range("D11").Value = "'0.07082"
$ws.Range("E11").Value = "'-0.66%"
$ws.Range("D12").Value = "'0.03205"
$ws.Range("E12").Value = "'0.41%"
$ws.Range("E13").Value = "'0.49%"
$ws.Range("D14").Value = "'0.001541"
$ws.Range("E14").Value = "'0.26%"
$ws.Range("D15").Value = "'0.0006045"
$ws.Range("E15").Value = "'-0.74%"
$ws.Range("D16").Value = "'0.006051"
$ws.Range("E16").Value = "'0.85%"
$ws.Range("D17").Value = "'3.516"
$ws.Range("E17").Value = "'0.59%"
$ws.Range("D18").Value = "'3.190"
$ws.Range("E18").Value = "'-0.75%"
$ws.Range("D19").Value = "'2.201"
$ws.Range("E19").Value = "'-1.08%"
$ws.Range("D20").Value = "'0.3072"
$ws.Range("E20").Value = "'-3.17%"
$ws.Range("D21").Value = "'0.1274"
$ws.Range("E21").Value = "'-2.60%"
$ws.Range("D22").Value = "'3.848"
$ws.Range("E22").Value = "'8.91%"
$ws.Range("D23").Value = "'0.04224"
$ws.Range("E23").Value = "'0.77%"
$ws.Range("D24").Value = "'0.001220"
$ws.Range("E24").Value = "'-0.54%"
$ws.Range("D25").Value = "'0.004279"
$ws.Range("E25").Value = "'-5.97%"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.10%"
$ws.Range("D27").Value = "'0.0001936"
$ws.Range("E27").Value = "'32.00%"
$ws.Range("D40").Value = "'0.03833"
$ws.Range("E40").Value = "'0.39%"
$ws.Range("D41").Value = "'0.006246"
$ws.Range("E41").Value = "'60.75%"
$ws.Range("D42").Value = "'0.1101"
$ws.Range("E42").Value = "'-0.04%"
$ws.Range("D43").Value = "'0.002302"
$ws.Range("E43").Value = "'-1.83%"
$ws.Range("D44").Value = "'0.01134"
$ws.Range("E44").Value = "'16.69%"
$ws.Range("D45").Value = "'0.00005458"
$ws.Range("E45").Value = "'0.90%"
$ws.Range("E46").Value = "'-0.08%"
$ws.Range("D47").Value = "'0.07775"
$ws.Range("E47").Value = "'-13.62%"
$ws.Range("E48").Value = "'6.88%"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.08%"
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.08%"
